$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 2 (originally "1：高维讲解。"): fix the typo first (讲解 -> 降解), which
# also splits the second run into "：" / "高维降解" / "。", then relabel "1" -> "A".
$tr.Replace("高维讲解", "高维降解") | Out-Null
$tr.Paragraphs(2).Runs(1).Text = "A"

# Paragraphs 3-6: simply relabel the leading number with a letter.
$tr.Paragraphs(3).Runs(1).Text = "B"
$tr.Paragraphs(4).Runs(1).Text = "C"
$tr.Paragraphs(5).Runs(1).Text = "D"
$tr.Paragraphs(6).Runs(1).Text = "E"
